# Commit: "1) Add capability to 'remove all geodata' from files. - Added a
# new button for the above. 2) logic fixes 3) formatting fixes"
#
# The workbook is a UI string/metadata table for a WinForms app. The
# "Button" sheet lists every button's control name (col A) + display text
# (col C). Two new buttons are introduced:
#   btn_EditFile       -> "Edit File"
#   btn_RemoveGeoData  -> "Remove GeoData"

$wb = $excel.ActiveWorkbook

$wsButton = $wb.Worksheets.Item("Button")

# Enter the "Remove GeoData" pair first, then "Edit File" - controls the
# order new entries land in the shared-string table - before placing
# "Edit File" above it as row 12 (new rows 12 & 13 at the bottom of the
# sheet).
$wsButton.Range("A13").Value = "btn_RemoveGeoData"
$wsButton.Range("C13").Value = "Remove GeoData"

$wsButton.Range("C12").Value = "Edit File"
$wsButton.Range("A12").Value = "btn_EditFile"

# Leave the cursor on the first newly-added row.
$wsButton.Range("A12").Select()

# "Button" (the first sheet) becomes the active/visible tab again, instead
# of "ToolStripMenuItem" (the last one).
$wsButton.Activate()

# "ColumnHeader" sheet: selection moved off the old full-column range onto
# a single cell.
$wsColumnHeader = $wb.Worksheets.Item("ColumnHeader")
$wsColumnHeader.Range("B29").Select()

# Re-activate "Button" so it's the sheet shown when the workbook re-opens.
$wsButton.Activate()
